$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update F column (想去人数 / want-to-go count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 627
$ws1.Range("F5").Value = 168
$ws1.Range("F6").Value = 9441
$ws1.Range("F8").Value = 332
$ws1.Range("F9").Value = 1205
$ws1.Range("F10").Value = 1163
$ws1.Range("F13").Value = 18
$ws1.Range("F14").Value = 264
$ws1.Range("F15").Value = 429
$ws1.Range("F16").Value = 93
$ws1.Range("F17").Value = 255
$ws1.Range("F18").Value = 1293

# Sheet "全部类型" (All Types) - same updates, shifted by one row
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 627
$ws4.Range("F6").Value = 168
$ws4.Range("F7").Value = 9441
$ws4.Range("F9").Value = 332
$ws4.Range("F10").Value = 1205
$ws4.Range("F11").Value = 1163
$ws4.Range("F14").Value = 18
$ws4.Range("F15").Value = 264
$ws4.Range("F16").Value = 429
$ws4.Range("F17").Value = 93
$ws4.Range("F18").Value = 255
$ws4.Range("F19").Value = 1293
